$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3575
$ws1.Range("F5").Value = 2213
$ws1.Range("F9").Value = 76
$ws1.Range("F10").Value = 64
$ws1.Range("F11").Value = 1316
$ws1.Range("F12").Value = 237
$ws1.Range("F13").Value = 1863
$ws1.Range("F14").Value = 137

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 3575
$ws4.Range("F5").Value = 2213
$ws4.Range("F10").Value = 76
$ws4.Range("F11").Value = 64
$ws4.Range("F14").Value = 1316
$ws4.Range("F15").Value = 237
$ws4.Range("F16").Value = 1863
$ws4.Range("F17").Value = 137
